$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.659.03"
$ws.Range("E2").Value = "  -1.13%  "

# Row 3
$ws.Range("D3").Value = "3.862.97"
$ws.Range("E3").Value = "  -2.08%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "'523.10"
$ws.Range("E5").Value = "  +6.50%  "

# Row 6
$ws.Range("D6").Value = "'140.86"
$ws.Range("E6").Value = "  -4.27%  "

# Row 7
$ws.Range("E7").Value = "  -2.70%  "

# Row 8
$ws.Range("E8").Value = "  +0.17%  "

# Row 9
$ws.Range("E9").Value = "  -3.45%  "

# Row 10
$ws.Range("D10").Value = "'0.167"
$ws.Range("E10").Value = "  -6.04%  "

# Row 11
$ws.Range("E11").Value = "  -8.17%  "

# Row 12
$ws.Range("D12").Value = "'41.59"
$ws.Range("E12").Value = "  -3.60%  "

# Row 13
$ws.Range("D13").Value = "'10.39"
$ws.Range("E13").Value = "  -1.02%  "

# Row 14
$ws.Range("D14").Value = "4.476.49"
$ws.Range("E14").Value = "  -2.07%  "

# Row 15
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "3.889.16"
$ws.Range("E15").Value = "  -1.47%  "

# Row 16
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'21.19"
$ws.Range("E16").Value = "  +6.44%  "

# Row 17
$ws.Range("D17").Value = "'14.12"
$ws.Range("E17").Value = "  -0.92%  "

# Row 18
$ws.Range("E18").Value = "  -2.16%  "

# Row 19
$ws.Range("E19").Value = "  +2.30%  "

# Row 20
$ws.Range("D20").Value = "68.642.41"
$ws.Range("E20").Value = "  -1.25%  "

# Row 21
$ws.Range("D21").Value = "'416.17"
$ws.Range("E21").Value = "  -5.28%  "

# Row 22
$ws.Range("D22").Value = "'3.52"
$ws.Range("E22").Value = "  +1.80%  "

# Row 23
$ws.Range("D23").Value = "'13.96"
$ws.Range("E23").Value = "  -4.25%  "

# Row 24
$ws.Range("D24").Value = "'86.77"
$ws.Range("E24").Value = "  -3.05%  "

# Row 25
$ws.Range("D25").Value = "'3.98"
$ws.Range("E25").Value = "  +6.31%  "

# Row 26
$ws.Range("D26").Value = "'11.51"
$ws.Range("E26").Value = "  -4.86%  "

# Row 27
$ws.Range("D27").Value = "'10.52"
$ws.Range("E27").Value = "  -5.93%  "

# Row 28
$ws.Range("D28").Value = "'35.55"
$ws.Range("E28").Value = "  -4.86%  "

# Row 29
$ws.Range("D29").Value = "'13.34"
$ws.Range("E29").Value = "  -1.37%  "

# Row 30
$ws.Range("D30").Value = "'679.17"
$ws.Range("E30").Value = "  -3.99%  "

# Row 31
$ws.Range("E31").Value = "  -5.21%  "

# Row 32
$ws.Range("D32").Value = "'6.86"
$ws.Range("E32").Value = "  +13.24%  "

# Row 33
$ws.Range("D33").Value = "'2.79"
$ws.Range("E33").Value = "  -3.47%  "

# Row 34
$ws.Range("D34").Value = "'67.20"
$ws.Range("E34").Value = "  +8.43%  "

# Row 35
$ws.Range("D35").Value = "'0.451"
$ws.Range("E35").Value = "  -4.18%  "

# Row 36
$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").Value = "'39.57"
$ws.Range("E36").Value = "  -3.00%  "

# Row 37
$ws.Range("B37").Value = "ThetaToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D37").Value = "'3.50"
$ws.Range("E37").Value = "  +14.06%  "

# Row 38
$ws.Range("D38").Value = "0.0₃0832"
$ws.Range("E38").Value = "  -8.23%  "

# Row 39
$ws.Range("E39").Value = "  -1.01%  "

# Row 40
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  +0.24%  "

# Row 41
$ws.Range("E41").Value = "  -0.17%  "

# Row 42
$ws.Range("D42").Value = "'0.0475"
$ws.Range("E42").Value = "  -3.12%  "

# Row 43
$ws.Range("D43").Value = "'3.14"
$ws.Range("E43").Value = "  +3.99%  "

# Row 44
$ws.Range("D44").Value = "'2.82"
$ws.Range("E44").Value = "  -4.48%  "

# Row 45
$ws.Range("D45").Value = "'3.40"
$ws.Range("E45").Value = "  +1.28%  "

# Row 46
$ws.Range("D46").Value = "'0.141"
$ws.Range("E46").Value = "  -1.44%  "

# Row 47
$ws.Range("D47").Value = "'3.00"
$ws.Range("E47").Value = "  -2.35%  "

# Row 48
$ws.Range("D48").Value = "'0.000274"
$ws.Range("E48").Value = "  +13.84%  "

# Row 49
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'145.28"
$ws.Range("E49").Value = "  +1.07%  "

# Row 50
$ws.Range("B50").Value = "LidoDAOToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D50").Value = "'3.28"
$ws.Range("E50").Value = "  -3.02%  "

# Row 51
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0340"
$ws.Range("E51").Value = "  -6.53%  "
